# Add "OpenFAST" / "turbsim" executable-path rows to the config sheet
# (new rows 18 and 19 right after the existing "OutSensors" row 17).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config")

$ws.Range("A18").Value = "OpenFAST"
$ws.Range("B18").Value = "openfast"
$ws.Range("A19").Value = "turbsim"
$ws.Range("B19").Value = "turbsim"

# Match the bold label styling used by the other label cells in column A
# (e.g. A17 "OutSensors") and the row heights of the surrounding rows.
$ws.Range("A18").Font.Bold = $true
$ws.Range("A19").Font.Bold = $true
$ws.Rows.Item(18).RowHeight = 13.8
$ws.Rows.Item(19).RowHeight = 13.8

# Make the config sheet the active sheet/tab with the new last row selected,
# mirroring the workbook-level activeTab switch captured in the edit.
$ws.Range("A19").Select()
